# Remove the "controleren van de configuratie op aanwezigheid van bekende
# kwetsbaarheden," bullet paragraph from the tools-list textbox on the
# "M16: Het project gebruikt tools voor vastgestelde taken" slide.
#
# Find it by content (rather than a hard-coded slide/shape/paragraph index)
# so the script is robust to minor structural differences, then delete the
# whole paragraph (its <a:p> node, including pPr/run) via TextRange.Delete().

$p = $ppt.ActivePresentation
$targetText = "controleren van de configuratie op aanwezigheid van bekende kwetsbaarheden,"

foreach ($s in $p.Slides) {
    foreach ($sh in $s.Shapes) {
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            $count = $tr.Paragraphs().Count
            # Walk backwards so deleting a paragraph doesn't shift the
            # indices of paragraphs we still need to inspect.
            for ($i = $count; $i -ge 1; $i--) {
                $para = $tr.Paragraphs($i, 1)
                # Each paragraph's .Text carries a trailing paragraph-mark
                # character (CR), so trim it before comparing.
                $txt = $para.Text.TrimEnd("`r", "`n", "`v")
                if ($txt -eq $targetText) {
                    $para.Delete()
                }
            }
        }
    }
}
